# Apply the "Add files via upload" edit to clothes.xlsx:
#  - two new product-type columns (AA: typ-upper, AB: typ-dress) with per-row
#    percentage values, mirroring the existing sex-male/sex-female columns
#  - D2/D3 (sex-male) bumped to 100 and E2/E3 (sex-female) reset to 0
#  - selection moved to AB12 (the one row where typ-dress wins)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: new columns AA (typ-upper) / AB (typ-dress) ---------------
$ws.Range("AA1").Value = "typ-upper"
$ws.Range("AB1").Value = "typ-dress"

# Match the vertical-center alignment style used by the rest of row 1, and
# also touch AC1 (empty, same style) so the sheet's used range/dimension
# extends out to column AC like the target workbook.
$ws.Range("AA1:AC1").VerticalAlignment = -4108

# --- corrected sex split for products 1 & 2 (rows 2-3) ---------------------
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 0

# --- new AA/AB (typ-upper / typ-dress) values for every data row -----------
$typeValues = @{
    2  = @(100, 0)
    3  = @(100, 0)
    4  = @(100, 0)
    5  = @(100, 0)
    6  = @(100, 0)
    7  = @(100, 0)
    8  = @(100, 0)
    9  = @(100, 0)
    10 = @(100, 0)
    11 = @(100, 0)
    12 = @(0, 100)
    13 = @(100, 0)
    14 = @(100, 0)
    15 = @(100, 0)
    16 = @(100, 0)
    17 = @(100, 0)
    18 = @(100, 0)
    19 = @(100, 0)
    20 = @(100, 0)
    21 = @(100, 0)
}

foreach ($row in $typeValues.Keys) {
    $vals = $typeValues[$row]
    $ws.Cells.Item($row, 27).Value = $vals[0]
    $ws.Cells.Item($row, 28).Value = $vals[1]
}

# --- view state: move the selection to AB12 --------------------------------
$ws.Range("AB12").Select()
